$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master")

# --- Row 2: sourcenum --- (B2=1, C2=2 unchanged; add D2=3)
$ws.Range("D2").Value = 3

# --- Row 3: type --- (B3 gaussian -> point; C3 stays point; add D3 = point)
$ws.Range("B3").Value = "point"
$ws.Range("C3").Value = "point"
$ws.Range("D3").Value = "point"

# --- Row 4: x pos [asec] --- (B4=0 unchanged; C4 10 -> 5; add D4=2)
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 2

# --- Row 5: y pos [asec] --- (B5=0 unchanged; C5 10 -> 0; add D5=-2)
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = -2

# --- Row 6: temperature --- (B6=60, C6=60 unchanged; add D6=50)
$ws.Range("D6").Value = 50

# --- Row 7: xwidth --- (B7=2 unchanged; add C7=3, D7=1)
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 1

# --- Row 8: ywidth --- (B8=2 unchanged; add C8=3, D8=1)
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 1

# --- Row 9: emissivity --- (B9=1 unchanged; C9 0 -> 1; add D9=1)
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1

# --- Row 10: linefreq --- (B10 30 -> 40; add C10=40; no D10)
$ws.Range("B10").Value = 40
$ws.Range("C10").Value = 40

# --- Row 11: spectrum --- (B11, C11 stay blackbody; add D11 = blackbody)
$ws.Range("D11").Value = "blackbody"

# Reselect to match the saved cursor position
$null = $ws.Range("D9").Select()
